$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Generated on" timestamp text in A2
$ws.Range("A2").Value = "Generated on: Fri Mar 26 23:51:27 WEST 2021"

# 2. Clear the sample data values in rows 5 and 6 (A:D), keeping A/B formatting
$ws.Range("A5:D6").ClearContents()

# 3. Remove the number-format styling from the Currency/Percent columns
#    (C and D) across rows 5-7, reverting them to the default style
$ws.Range("C5:D7").ClearFormats()
